# Remove the three leading "site navigation" paragraphs (Home / Back to
# Home / Download Word Document hyperlinks) that preceded the
# "Table of Contents" heading. Everything else in the body is unchanged;
# the bookmark ids shift down automatically as a consequence of removing
# the bookmark-less content ahead of the "table-of-contents" bookmark.

$d = $word.ActiveDocument

$first = $d.Paragraphs.Item(1)
$last  = $d.Paragraphs.Item(3)

# Sanity-check we are about to remove the expected paragraphs before
# deleting anything.
if ($first.Range.Text -notmatch "Home" -or $last.Range.Text -notmatch "Download Word Document") {
    throw "Unexpected document content; aborting to avoid deleting the wrong paragraphs."
}

$rng = $d.Range($first.Range.Start, $last.Range.End)
$rng.Delete()
